# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.476.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.13%  "
$ws.Range("D3").Value = "'2.065.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'252.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("D6").Value = "'0.651"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("D7").Value = "'65.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.47%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.399"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.86%  "
$ws.Range("D10").Value = "'59.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").Value = "'0.0824"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.90%  "
$ws.Range("D12").Value = "'0.105"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "'0.925"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("D14").Value = "'23.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +23.13%  "
$ws.Range("D15").Value = "'14.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "'2.351.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").Value = "'5.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("D18").Value = "'2.042.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("D19").Value = "'37.391.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.15%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'73.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0₃0913"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.14%  "
$ws.Range("D22").Value = "'5.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.29%  "
$ws.Range("D23").Value = "'239.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "'2.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.10%  "
$ws.Range("D27").Value = "'10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.98%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.03%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'162.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").Value = "'0.128"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +29.19%  "
$ws.Range("D31").Value = "'0.123"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("D32").Value = "'5.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").Value = "'1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.90%  "
$ws.Range("D34").Value = "'0.0632"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("D35").Value = "'4.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.80%  "
$ws.Range("D36").Value = "'2.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'6.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.39%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +3.25%  "
$ws.Range("D40").Value = "'3.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +32.20%  "
$ws.Range("E41").Value = "  +4.61%  "
$ws.Range("E42").Value = "  +5.65%  "
$ws.Range("E43").Value = "  +4.94%  "
$ws.Range("D44").Value = "'1.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.66%  "
$ws.Range("D45").Value = "'17.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.74%  "
$ws.Range("D46").Value = "'0.0220"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D47").Value = "'95.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("D48").Value = "'7.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").Value = "'1.395.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").Value = "'2.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "'46.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
